$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fix header text: "preguntas" -> "pregunta" (singular)
$ws.Range("B1").Value = "pregunta"

# Move active selection from A2 to B2
$ws.Range("B2").Select()
